$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "VDFUP001"
$ws.Range("C2").Value = "POLVOS - FABRICACIONES DEL 04 AL 08 DE SEP.xlsx"

# Remove rows 3 through 13 (old extra error rows)
$ws.Rows("3:13").Delete()

# Adjust column widths: A -> 10, C -> 49 (accounting for Excel's ~0.83 padding offset)
$ws.Columns("A").ColumnWidth = 9.17
$ws.Columns("C").ColumnWidth = 48.17
